# Updated cryptos list on Thu May  9 11:36:33 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# coins that moved, plus a rank swap between PancakeSwap and FirstDigitalUSD
# (rows 29/30 traded places). All of these columns are stored as text in the
# sheet (e.g. "60.974.15" isn't a real number, and the percentages keep their
# padding spaces), so a leading apostrophe is used wherever the new value
# would otherwise look like a plain number to Excel - this keeps those cells
# as text instead of letting them get auto-converted to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.985.12"
$ws.Range("D3").Value = "2.973.09"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'595.21"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'142.17"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").Value = "2.971.30"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").Value = "'5.96"
$ws.Range("E11").Value = "  +4.45%  "
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "'33.90"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D16").Value = "3.463.54"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "61.124.06"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "'6.82"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").Value = "2.970.02"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'449.20"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'7.26"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "'81.99"
$ws.Range("D25").Value = "'2.16"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").Value = "'10.36"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").Value = "'11.87"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.66"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").Value = "'27.24"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "0.0₃0801"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "'49.93"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'2.03"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +6.28%  "
$ws.Range("D43").Value = "'385.29"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "'0.0347"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("D46").Value = "'38.28"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "2.690.72"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").Value = "'129.88"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.11%  "
